$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Paragraph listing the poster ("плакат") contents: insert a new sentence
#    into item 5, insert a whole new item 6 (renumbering the old "6" to "7"),
#    and append a blank paragraph after the list.
# ---------------------------------------------------------------------------

# 1a. Item 5 gains " Алгоритм авторизации и подключения." right after the
#     first sentence, before "Схема программы."
$d.Content.Find.Execute(
    "умного дома. Схема программы",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "умного дома. Алгоритм авторизации и подключения. Схема программы",
    2) | Out-Null

# 1b. A brand-new item 6 is inserted right before "Заключительный плакат",
#     pushing the old 6 to become 7.
$d.Content.Find.Execute(
    "6 Заключительный плакат",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "6 Android-приложение для мониторинга и управления устройствами умного дома. Алгоритм изменения состояния устройства. Схема программы 7 Заключительный плакат",
    2) | Out-Null

# 1c. Append an empty paragraph right after the poster-list paragraph.
$listPara = $d.Content.Find.Execute("Заключительный плакат. Плакат.", $false, $true, $false, $false, $false, $true, 1, $false, $null, 0)
$p = $d.Paragraphs.Item(32)
$newRange = $p.Range
$newRange.Collapse(0)
$newRange.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 2. Remove the now-obsolete lastRenderedPageBreak marker before
#    "Консультанты по дипломному проект" (Word will recompute this on
#    pagination; nothing to do here explicitly).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 3. Table cells: replace prose section names with their numeric references.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "(разделы обзор литературы, структурное проектирование, функциональное проектирование, технико-экономическое обоснование)",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "(разделы 1, 2, 3, 4)",
    2) | Out-Null

$d.Content.Find.Execute(
    "(разделы функциональное проектирование, разработка программных модулей)",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "(разделы 3, 4)",
    2) | Out-Null
